# Update cryptocurrency price/volume data (refresh from source)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.447.62'
$ws.Range("E2").Value = '  -0.97%  '
$ws.Range("D3").Value = '2.781.83'
$ws.Range("E3").Value = '  -0.31%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '352.64'
$ws.Range("E5").Value = '  -2.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '107.96'
$ws.Range("E6").Value = '  -1.54%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.550'
$ws.Range("E7").Value = '  -1.56%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.620'
$ws.Range("E9").Value = '  +4.98%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.21'
$ws.Range("E10").Value = '  -2.41%  '
$ws.Range("E11").Value = '  +1.37%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0834'
$ws.Range("E12").Value = '  -1.80%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.87'
$ws.Range("E13").Value = '  +1.80%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.74'
$ws.Range("E14").Value = '  +2.33%  '
$ws.Range("D15").Value = '3.217.16'
$ws.Range("E15").Value = '  -0.33%  '
$ws.Range("D16").Value = '2.794.42'
$ws.Range("E16").Value = '  -0.51%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.926'
$ws.Range("E17").Value = '  -1.39%  '
$ws.Range("D18").Value = '51.441.31'
$ws.Range("E18").Value = '  -0.91%  '
$ws.Range("E19").Value = '  +3.29%  '
$ws.Range("E20").Value = '  +0.74%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.34'
$ws.Range("E21").Value = '  +1.55%  '
$ws.Range("D22").Value = '0.0₃0967'
$ws.Range("E22").Value = '  -0.87%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.51'
$ws.Range("E23").Value = '  +0.12%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '266.36'
$ws.Range("E24").Value = '  -1.28%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.76'
$ws.Range("E25").Value = '  +0.04%  '
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.85'
$ws.Range("E28").Value = '  +2.00%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.26'
$ws.Range("E29").Value = '  -0.44%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '37.01'
$ws.Range("E30").Value = '  +7.30%  '
$ws.Range("E31").Value = '  -2.51%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.25'
$ws.Range("E32").Value = '  +8.96%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '51.81'
$ws.Range("E33").Value = '  -0.33%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.67'
$ws.Range("E34").Value = '  +8.06%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0441'
$ws.Range("E35").Value = '  -6.26%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0848'
$ws.Range("E36").Value = '  +0.20%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  +0.00%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.48'
$ws.Range("E38").Value = '  -3.03%  '
$ws.Range("E39").Value = '  -3.13%  '
$ws.Range("E40").Value = '  -1.87%  '
$ws.Range("E42").Value = '  -5.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '120.16'
$ws.Range("E43").Value = '  +0.57%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '21.90'
$ws.Range("E44").Value = '  -0.58%  '
$ws.Range("E45").Value = '  -2.74%  '
$ws.Range("D46").Value = '2.129.30'
$ws.Range("E46").Value = '  +2.15%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.36'
$ws.Range("E47").Value = '  +3.33%  '
$ws.Range("E48").Value = '  +5.13%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.226'
$ws.Range("E49").Value = '  +18.31%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.47'
$ws.Range("E50").Value = '  -5.56%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.893'
$ws.Range("E51").Value = '  -6.82%  '
